$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.371.03"
$ws.Range("E2").Value = "  +0.42%  "
$ws.Range("D3").Value = "2.935.28"
$ws.Range("E3").Value = "  +0.34%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'594.95"
$ws.Range("E5").Value = "  +0.60%  "
$ws.Range("D6").Value = "'145.06"
$ws.Range("E6").Value = "  -0.17%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  -0.35%  "
$ws.Range("D9").Value = "'6.99"
$ws.Range("E9").Value = "  +2.25%  "
$ws.Range("E10").Value = "  -0.90%  "
$ws.Range("D11").Value = "'0.441"
$ws.Range("E11").Value = "  -0.44%  "
$ws.Range("E12").Value = "  -0.42%  "
$ws.Range("D13").Value = "'33.67"
$ws.Range("E13").Value = "  -0.24%  "
$ws.Range("E14").Value = "  +0.57%  "
$ws.Range("D15").Value = "3.454.98"
$ws.Range("E15").Value = "  +1.17%  "
$ws.Range("D16").Value = "61.374.86"
$ws.Range("E16").Value = "  +0.49%  "
$ws.Range("D17").Value = "'6.74"
$ws.Range("E17").Value = "  +0.18%  "
$ws.Range("D18").Value = "2.936.31"
$ws.Range("E18").Value = "  +0.15%  "
$ws.Range("D19").Value = "'433.59"
$ws.Range("E19").Value = "  +0.50%  "
$ws.Range("D20").Value = "'13.52"
$ws.Range("E20").Value = "  +0.12%  "
$ws.Range("D21").Value = "'0.681"
$ws.Range("E21").Value = "  -0.36%  "
$ws.Range("D22").Value = "'7.17"
$ws.Range("E22").Value = "  +1.03%  "
$ws.Range("D23").Value = "'81.96"
$ws.Range("E23").Value = "  +1.29%  "
$ws.Range("D24").Value = "'10.95"
$ws.Range("E24").Value = "  +1.10%  "
$ws.Range("D25").Value = "'2.22"
$ws.Range("E25").Value = "  -0.18%  "
$ws.Range("D26").Value = "'11.88"
$ws.Range("E26").Value = "  -2.37%  "
$ws.Range("E27").Value = "  -0.16%  "
$ws.Range("E28").Value = "  -3.99%  "
$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("E30").Value = "  -1.98%  "
$ws.Range("E31").Value = "  +2.95%  "
$ws.Range("D32").Value = "'26.77"
$ws.Range("E32").Value = "  +0.69%  "
$ws.Range("D33").Value = "'1.00"
$ws.Range("E33").Value = "  -0.10%  "
$ws.Range("D34").Value = "0.0₃0887"
$ws.Range("E34").Value = "  +2.40%  "
$ws.Range("D35").Value = "'1.01"
$ws.Range("E35").Value = "  -0.15%  "
$ws.Range("D36").Value = "'5.67"
$ws.Range("E36").Value = "  +0.80%  "
$ws.Range("D37").Value = "'3.02"
$ws.Range("E37").Value = "  -2.80%  "
$ws.Range("D38").Value = "'2.03"
$ws.Range("E38").Value = "  +0.72%  "
$ws.Range("E39").Value = "  -0.58%  "
$ws.Range("D40").Value = "'8.64"
$ws.Range("E40").Value = "  +0.10%  "
$ws.Range("E41").Value = "  +6.82%  "
$ws.Range("D42").Value = "'0.284"
$ws.Range("E42").Value = "  -1.62%  "
$ws.Range("D43").Value = "'0.0349"
$ws.Range("E43").Value = "  -0.12%  "
$ws.Range("D44").Value = "'372.92"
$ws.Range("E44").Value = "  -2.60%  "
$ws.Range("D45").Value = "2.716.36"
$ws.Range("E45").Value = "  +0.30%  "
$ws.Range("D46").Value = "'133.48"
$ws.Range("E46").Value = "  +2.67%  "
$ws.Range("D48").Value = "'24.00"
$ws.Range("E48").Value = "  -1.52%  "
$ws.Range("D49").Value = "'0.105"
$ws.Range("E49").Value = "  -1.23%  "
$ws.Range("E50").Value = "  -0.67%  "
$ws.Range("D51").Value = "'0.125"
$ws.Range("E51").Value = "  -0.65%  "
